$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 79, shifting existing rows 79-108 down to 80-109.
$ws.Rows.Item(79).Insert()

# Copy formatting (esp. the date-format style on column D) from the row
# immediately below (the original row 79, now row 80) into the new row 79.
$ws.Range("A80:R80").Copy()
$ws.Range("A79:R79").PasteSpecial(-4122) | Out-Null

# Fill the new row 79 with the data for the new weekly record.
$ws.Cells.Item(79, 1).Value = 9
$ws.Cells.Item(79, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(79, 3).Value = "Metropolitana"
$ws.Cells.Item(79, 4).Value = 45119
$ws.Cells.Item(79, 5).Value = 13
$ws.Cells.Item(79, 6).Value = 100112029
$ws.Cells.Item(79, 7).Value = "Orégano"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 16
$ws.Cells.Item(79, 11).Value = 21000
$ws.Cells.Item(79, 12).Value = 21000
$ws.Cells.Item(79, 13).Value = 21000
$ws.Cells.Item(79, 14).Value = "$/docena de atados"
$ws.Cells.Item(79, 15).Value = "Región Metropolitana"
$ws.Cells.Item(79, 16).Value = 7000
$ws.Cells.Item(79, 17).Value = 3
$ws.Cells.Item(79, 18).Value = "Hortaliza"
